$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$expected = @(
  "1+50=",
  "66-19=",
  "60+18=",
  "44-44=",
  "93-22=",
  "8+61=",
  "7+11=",
  "74-0=",
  "37+22=",
  "19+37=",
  "61-44=",
  "72-37=",
  "94-72=",
  "47-16=",
  "19+27=",
  "82-40=",
  "47-7=",
  "76+22=",
  "0+44=",
  "71-46=",
  "36-35=",
  "43+43=",
  "98-60=",
  "79-40=",
  "33+1=",
  "96-37=",
  "83-43=",
  "79-59=",
  "78-25=",
  "30+49=",
  "20+22=",
  "81-57=",
  "51+21=",
  "54-23=",
  "4+66=",
  "36-16=",
  "89-4=",
  "96-5=",
  "34+52=",
  "14+26=",
  "76-28=",
  "24+58=",
  "41-17=",
  "77-49=",
  "41+38=",
  "37+10=",
  "38-33=",
  "96-32=",
  "75+2=",
  "1+12=",
  "23-1=",
  "76-68=",
  "6+32=",
  "69+28=",
  "76-52=",
  "48+35=",
  "40+51=",
  "31+61=",
  "91-19=",
  "61-32=",
  "88-44=",
  "54-2=",
  "6+24=",
  "91-32=",
  "41+20=",
  "7+59=",
  "5+28=",
  "89+2=",
  "21-13=",
  "66-25=",
  "78-17=",
  "96-24=",
  "75+8=",
  "47+19=",
  "90-33=",
  "0+34=",
  "74+18=",
  "16+24=",
  "75+20=",
  "17+65=",
  "94-83=",
  "48-43=",
  "75-17=",
  "32-31=",
  "8-6=",
  "92-92=",
  "25+40=",
  "16+38=",
  "84-9=",
  "37+21=",
  "65-1=",
  "13+84=",
  "77-7=",
  "88-65=",
  "72-63=",
  "2+8=",
  "23+14=",
  "95-65=",
  "96-49=",
  "31+67="
)
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
  for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $expected[$idx]
    $idx++
  }
}
Write-Host "Done: applied " $idx " replacements"
